$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 5): insert "Start Uur" / "Stop Uur" columns before the
# existing "Gespendeerde Uren" / "Notes" columns, shifting those two right.
# Use Copy() (value+format) rather than reading .Value back, since .Value
# reads don't resolve to plain scalars in this host.
$ws.Range("E5").Copy($ws.Range("G5"))
$ws.Range("D5").Copy($ws.Range("F5"))
$ws.Range("C5").Copy($ws.Range("E5"))
$ws.Range("D5").Value = "Start Uur"
$ws.Range("E5").Value = "Stop Uur"

# --- Row 6: journal entry #1 ---
$ws.Range("C6").Value = "Opstarten research project"
$ws.Range("D6").Value = 0.5493055555555556
$ws.Range("E6").Value = 0.61388888888888882
$ws.Range("F6").Formula = "=E6-D6"
$ws.Range("G6").Value = "Besloten om python research te doen"

# --- Row 7-11: journal entry #2, with a filled-down (shared) formula in F,
# afterwards blanked out in F8:F11 but keeping the number format. ---
$ws.Range("D7").Value = 0.6694444444444444
$ws.Range("E7").Value = 0.71736111111111101
$ws.Range("G7").Value = "Maken mappenstructuur en GIT"
$ws.Range("F7:F11").Formula = "=E7-D7"
$ws.Range("F8:F11").ClearContents()

# --- Row 8 ---
$ws.Range("G8").Value = "Kiezen leermethode/platform: ""Codecademy"""

# --- Row 9 ---
$ws.Range("G9").Value = "Afgewerkte Codecademy lessen:"

# --- Apply the new time format (h:mm) to the Start/Stop Uur + Gespendeerde
# Uren (now column F) ranges, and the lone placeholder far below. ---
$ws.Range("D6:E8").NumberFormat = "h:mm"
$ws.Range("F6:F14").NumberFormat = "h:mm"
$ws.Range("F28").NumberFormat = "h:mm"

# --- Column widths (best achievable match; this host only supports the
# 1/6-character granularity behind ColumnWidth). ---
$ws.Columns.Item(4).ColumnWidth = 7.833333333333333
$ws.Columns.Item(5).ColumnWidth = 7.666666666666667
$ws.Columns.Item(6).ColumnWidth = 18.333333333333332
$ws.Columns.Item(7).ColumnWidth = 155.5

# --- Selection ---
$ws.Range("G12").Select()
